# Word Online couldn't parse numPr correctly when <w:numId> was emitted
# before <w:ilvl>. Re-apply each list paragraph's own current list level
# (a value no-op) so the document is re-serialized with <w:ilvl> before
# <w:numId> inside every <w:numPr>, matching the fixed element order.
$d = $word.ActiveDocument

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.ListFormat.ListType -ne 0) {
        $lvl = $p.Range.ListFormat.ListLevelNumber
        $p.Range.ListFormat.ListLevelNumber = $lvl
    }
}
